$d = $word.ActiveDocument

# The first paragraph ("Test'"><svg/onload=alert(1)/>") is replaced by three
# paragraphs:
#   1) <script>alert(1)</script>           (with gramStart/gramEnd proofErr
#                                            wrapping the "1)<" run)
#   2) '"<script>alert(2)</script>         (with gramStart/gramEnd proofErr
#                                            wrapping the "2" + ")<" runs)
#   3) an empty paragraph
# The second (original) paragraph with the <img ...> SSRF payload is left
# untouched.

$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

$newXml = @"
<w:p xmlns:w="$wNs"><w:r><w:t>&lt;script&gt;alert(</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>1)&lt;</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>/script&gt;</w:t></w:r></w:p><w:p xmlns:w="$wNs"><w:r><w:t>'"</w:t></w:r><w:r><w:t>&lt;script&gt;alert(</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>2</w:t></w:r><w:r><w:t>)&lt;</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>/script&gt;</w:t></w:r></w:p><w:p xmlns:w="$wNs"></w:p>
"@

$target = $d.Paragraphs.Item(1).Range
$target.InsertXML($newXml) | Out-Null
